$d = $word.ActiveDocument

# Ensure edits are applied as plain changes, not tracked revisions.
$d.TrackRevisions = $false

# 1) In the "MIND SATELLITE PREVENTION SECURITY SYSTEMS" heading (bold, 12pt run),
#    merge the "SATELLITE" run and the following space run into a single run's text.
#    This paragraph is the 15th paragraph in the document.
$p1 = $d.Paragraphs.Item(15).Range
$p1.Find.Execute("SATELLITE ", $true, $false, $false, $false, $false, $true, 1, $false, "SATELLITE ", 2) | Out-Null

# 2) Change the year reference from 2022 to 2023 in the following paragraph.
$p2 = $d.Paragraphs.Item(16).Range
$p2.Find.Execute("2022", $true, $false, $false, $false, $false, $true, 1, $false, "2023", 2) | Out-Null
